# Tire Type Filtering update — recomputed signal-processing pipeline values.
#
# Step1_Data  : raw per-segment signal values (Signal_Value_58..92 in cols B..AJ)
#               were recomputed by the updated audio->csv extraction script.
# Step2_Sj    : row-wise running cumulative sum of Step1_Data (derived).
# Step3_DataPts_0.5 / 0.7 / 0.8 / 0.9 :
#               for each intensity threshold, the first signal position whose
#               cumulative value (Step2_Sj) exceeds the threshold, the
#               cumulative value there, and the resulting pulse width.
#
# This script rewrites Step1_Data with the new raw values, then recomputes
# Step2_Sj and the four Step3_DataPts_* sheets from that new data so every
# derived sheet stays internally consistent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------
function ColLetterToIndex($letters) {
    $idx = 0
    for ($i = 0; $i -lt $letters.Length; $i++) {
        $ch = $letters.Substring($i, 1)
        $idx = $idx * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $idx
}

# ---------------------------------------------------------------------------
# New raw Step1_Data values (only the cells that actually changed)
# ---------------------------------------------------------------------------
$newStep1Values = @{}
$newStep1Values["E2"] = 0.1806491037374188
$newStep1Values["F2"] = 0.2919699818052948
$newStep1Values["G2"] = 0.008992064282946346
$newStep1Values["I2"] = 0.03178900251222747
$newStep1Values["J2"] = 0.01015438944110902
$newStep1Values["K2"] = 0.01484969935693541
$newStep1Values["M2"] = 0.07014322571309826
$newStep1Values["N2"] = 0.05654443852874948
$newStep1Values["Q2"] = 0.0007746774020523915
$newStep1Values["S2"] = 0.03231326094201436
$newStep1Values["T2"] = 0.06829717391734537
$newStep1Values["V2"] = 0.03170028378832641
$newStep1Values["X2"] = 0.009999374404804079
$newStep1Values["AC2"] = 0.0183097744726793
$newStep1Values["AD2"] = 0.09866332594815344
$newStep1Values["AE2"] = 0.0224382759445493
$newStep1Values["AI2"] = 0.05241194780229576
$newStep1Values["E3"] = 0.1152671548090895
$newStep1Values["G3"] = 0.2956481462747582
$newStep1Values["H3"] = 0.0005599140615966256
$newStep1Values["I3"] = 0.02069082663766454
$newStep1Values["L3"] = 0.02597367668428472
$newStep1Values["M3"] = 0.01843487522642335
$newStep1Values["N3"] = 0.01117012641684706
$newStep1Values["O3"] = 0.1689122403933956
$newStep1Values["T3"] = 0.06067764504195666
$newStep1Values["U3"] = 0.04190095387099034
$newStep1Values["W3"] = 0.007773049800268729
$newStep1Values["Y3"] = 0.008920858445347565
$newStep1Values["AD3"] = 0.04176022220991234
$newStep1Values["AE3"] = 0.1273595999570336
$newStep1Values["AI3"] = 0.005755767935039623
$newStep1Values["AJ3"] = 0.04919494223539166
$newStep1Values["D4"] = 0.03184889134399857
$newStep1Values["E4"] = 0.0660802745501909
$newStep1Values["F4"] = 0.2814520017784889
$newStep1Values["H4"] = 0.01678270386237836
$newStep1Values["I4"] = 0.05109687159581133
$newStep1Values["K4"] = 0.00336669520463865
$newStep1Values["L4"] = 0.125192493619958
$newStep1Values["M4"] = 0.005289037054195332
$newStep1Values["N4"] = 0.06432171563346049
$newStep1Values["O4"] = 0.001440970016561723
$newStep1Values["P4"] = 0.005800679683419603
$newStep1Values["S4"] = 0.08324430015123112
$newStep1Values["T4"] = 0.02621079021411546
$newStep1Values["U4"] = 0.0002956634246393913
$newStep1Values["V4"] = 0.01179594107374536
$newStep1Values["Z4"] = 0.005240669574220732
$newStep1Values["AA4"] = 0.0004118597115124203
$newStep1Values["AC4"] = 0.02774845040082538
$newStep1Values["AD4"] = 0.127950606889279
$newStep1Values["AE4"] = 0.0138311525017559
$newStep1Values["AG4"] = 0.01738256668182936
$newStep1Values["AI4"] = 0.03321566503374412
$newStep1Values["E5"] = 0.2765696903328075
$newStep1Values["F5"] = 0.1593966693624001
$newStep1Values["G5"] = 0.04635337625608153
$newStep1Values["I5"] = 0.03738264665979879
$newStep1Values["K5"] = 0.009822858848574034
$newStep1Values["L5"] = 0.0732565474801822
$newStep1Values["M5"] = 0.006670411058291809
$newStep1Values["N5"] = 0.03921662726288724
$newStep1Values["O5"] = 0.03816842203002564
$newStep1Values["P5"] = 0.002415493653806708
$newStep1Values["S5"] = 0.02624769953171851
$newStep1Values["T5"] = 0.06545724656504656
$newStep1Values["V5"] = 0.01821426470059954
$newStep1Values["X5"] = 0.00971833519164832
$newStep1Values["AA5"] = 0.005882868038165722
$newStep1Values["AC5"] = 0.001810024550478369
$newStep1Values["AD5"] = 0.09398392014267906
$newStep1Values["AE5"] = 0.05858755177557829
$newStep1Values["AI5"] = 0.03084534655922999
$newStep1Values["E6"] = 0.2289901104482266
$newStep1Values["F6"] = 0.03425835099463566
$newStep1Values["G6"] = 0.1448711024940938
$newStep1Values["I6"] = 0.01362512119418069
$newStep1Values["L6"] = 0.09131402427987571
$newStep1Values["M6"] = 0.01661562580095783
$newStep1Values["N6"] = 0.02385911854767411
$newStep1Values["O6"] = 0.1016553528192589
$newStep1Values["S6"] = 0.008606666903703198
$newStep1Values["T6"] = 0.07820614351570217
$newStep1Values["U6"] = 0.02381478051187402
$newStep1Values["V6"] = 0.007263687633896638
$newStep1Values["X6"] = 0.0002233115756684116
$newStep1Values["AA6"] = 0.0111689939891105
$newStep1Values["AD6"] = 0.06682411997880898
$newStep1Values["AE6"] = 0.11317262287386
$newStep1Values["AI6"] = 0.01478366052008975
$newStep1Values["AJ6"] = 0.02074720591838335

# ---------------------------------------------------------------------------
# 1) Write the new raw values into Step1_Data
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Step1_Data")

foreach ($key in $newStep1Values.Keys) {
    if ($key -match "^([A-Z]+)(\d+)$") {
        $colLetters = $matches[1]
        $rowNum = [int]$matches[2]
        $colIdx = ColLetterToIndex $colLetters
        $wsData.Cells.Item($rowNum, $colIdx).Value2 = $newStep1Values[$key]
    }
}

# ---------------------------------------------------------------------------
# 2) Recompute Step2_Sj = running cumulative sum of Step1_Data (cols B..AJ)
# ---------------------------------------------------------------------------
$wsSj = $wb.Worksheets.Item("Step2_Sj")

$firstCol = 2   # column B
$lastCol = 36   # column AJ

# store the recomputed cumulative rows so the Step3 sheets below can reuse
# them without re-reading from the sheet (keeps everything in one pass).
$cumRows = @{}

for ($row = 2; $row -le 6; $row++) {
    $running = 0
    $rowVals = @()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $v = $wsData.Cells.Item($row, $col).Value2
        $running = $running + $v
        $wsSj.Cells.Item($row, $col).Value2 = $running
        $rowVals += $running
    }
    $cumRows[$row] = $rowVals
}

# ---------------------------------------------------------------------------
# 3) Recompute the Step3_DataPts_* sheets (Point_Exceeds_Index / Cumulative
#    Value / Pulse_Width) for every intensity threshold sheet.
#    Column layout: A Segment_ID, B Intensity_Threshold,
#    C First_Noticeable_Increase_Index, D Point_Exceeds_Index,
#    E First_Noticeable_Increase_Cumulative_Value,
#    F Point_Exceeds_Cumulative_Value, G Pulse_Width, H.. unchanged metadata.
# ---------------------------------------------------------------------------
$step3Sheets = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

foreach ($sheetName in $step3Sheets) {
    $ws3 = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 6; $row++) {
        $threshold = $ws3.Cells.Item($row, 2).Value2          # B: Intensity_Threshold
        $firstIncreaseIdx = $ws3.Cells.Item($row, 3).Value2   # C: stays fixed

        $rowVals = $cumRows[$row]

        $exceedPos = -1
        $exceedVal = 0
        for ($i = 0; $i -lt $rowVals.Length; $i++) {
            if ($exceedPos -eq -1) {
                if ($rowVals[$i] -gt $threshold) {
                    $exceedPos = $i + 1   # 1-based position, B = position 1
                    $exceedVal = $rowVals[$i]
                }
            }
        }

        $ws3.Cells.Item($row, 4).Value2 = $exceedPos                      # D: Point_Exceeds_Index
        $ws3.Cells.Item($row, 6).Value2 = $exceedVal                      # F: Point_Exceeds_Cumulative_Value
        $ws3.Cells.Item($row, 7).Value2 = $exceedPos - $firstIncreaseIdx  # G: Pulse_Width
    }
}

Write-Output "Recomputed Step1_Data, Step2_Sj, and Step3_DataPts_* sheets."
